# Adiciona comarca aos documentos de prorrogações
#
# The date line at the top of the document used to read a hard-coded
# "Rio de Janeiro, {{ data_hoje }}." -- replace the hard-coded city name
# with a "{{ comarca }}" merge placeholder (like the rest of the
# template's merge fields) and highlight it in the same attention-red
# used for the other placeholder on that line, so authors notice it
# needs to be filled in / reviewed per case.

$d = $word.ActiveDocument

# --- 1. Swap the literal city name for the {{ comarca }} placeholder ---
$found = $d.Content.Find.Execute(
    "Rio de Janeiro", $true, $false, $false, $false, $false,
    $true, 1, $false, "{{ comarca }}", 2)

if (-not $found) {
    throw "Expected to find the literal 'Rio de Janeiro' run to replace with {{ comarca }}."
}

# --- 2. Recolor just that placeholder run to FF3333 (BGR 3355647) ---
# Re-locate the freshly-inserted placeholder text so the color change is
# scoped to exactly that run (matches the sibling ", {{ data_hoje }}" run
# which already uses a highlight color to flag a merge field).
$rngComarca = $d.Content
$foundAgain = $rngComarca.Find.Execute("{{ comarca }}")
if ($foundAgain) {
    $rngComarca.Font.Color = 3355647
}

Write-Output ("comarca placeholder inserted=" + $found + " recolored=" + $foundAgain)

# --- 3. Normal style's body-text color moves off "Automatic" onto an
#        explicit near-black (00000A), matching the rest of the template
#        styles which already pin an explicit color instead of "auto". ---
try {
    $normal = $d.Styles.Item("Normal")
    $normal.Font.Color = 655360   # wdColor BGR for RGB(00,00,0A) = 0x00000A
    Write-Output ("Normal style font color=" + $normal.Font.Color)
} catch {
    Write-Output ("Normal style color update skipped: " + $_)
}
